$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: convert existing rows 585-599 columns C and E from text to numeric ---
for ($r = 585; $r -le 599; $r++) {
    $c = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 3).Value = [double]$c
    $e = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value = [double]$e
}

# --- Part 2: append new rows 600-608 (Monday, May 17, 2021) with numeric scores ---
$numericRows = @(
    @(600, 'Monday, May 17, 2021', 'San Francisco Giants', 1, 'Cincinnati Reds', 0),
    @(601, 'Monday, May 17, 2021', 'New York Mets', 0, 'Atlanta Braves', 0),
    @(602, 'Monday, May 17, 2021', 'Washington Nationals', 0, 'Chicago Cubs', 1),
    @(603, 'Monday, May 17, 2021', 'Chicago White Sox', 3, 'Minnesota Twins', 0),
    @(604, 'Monday, May 17, 2021', 'New York Yankees', 1, 'Texas Rangers', 0),
    @(605, 'Monday, May 17, 2021', 'Cleveland Indians', 1, 'Los Angeles Angels', 1),
    @(606, 'Monday, May 17, 2021', 'Arizona Diamondbacks', 0, 'Los Angeles Dodgers', 0),
    @(607, 'Monday, May 17, 2021', 'Colorado Rockies', 0, 'San Diego Padres', 3),
    @(608, 'Monday, May 17, 2021', 'Detroit Tigers', 0, 'Seattle Mariners', 0)
)

foreach ($row in $numericRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# --- Part 3: append new rows 609-623 (Tuesday, May 18, 2021) with text ("inlineStr"-style) scores ---
$textRows = @(
    @(609, 'Tuesday, May 18, 2021', 'San Francisco Giants', 0, 'Cincinnati Reds', 1),
    @(610, 'Tuesday, May 18, 2021', 'Tampa Bay Rays', 0, 'Baltimore Orioles', 0),
    @(611, 'Tuesday, May 18, 2021', 'Miami Marlins', 0, 'Philadelphia Phillies', 0),
    @(612, 'Tuesday, May 18, 2021', 'New York Mets', 0, 'Atlanta Braves', 0),
    @(613, 'Tuesday, May 18, 2021', 'Boston Red Sox', 0, 'Toronto Blue Jays', 0),
    @(614, 'Tuesday, May 18, 2021', 'Washington Nationals', 0, 'Chicago Cubs', 0),
    @(615, 'Tuesday, May 18, 2021', 'Chicago White Sox', 1, 'Minnesota Twins', 0),
    @(616, 'Tuesday, May 18, 2021', 'Pittsburgh Pirates', 0, 'St. Louis Cardinals', 2),
    @(617, 'Tuesday, May 18, 2021', 'New York Yankees', 0, 'Texas Rangers', 0),
    @(618, 'Tuesday, May 18, 2021', 'Milwaukee Brewers', 0, 'Kansas City Royals', 0),
    @(619, 'Tuesday, May 18, 2021', 'Cleveland Indians', 5, 'Los Angeles Angels', 1),
    @(620, 'Tuesday, May 18, 2021', 'Houston Astros', 2, 'Oakland Athletics', 1),
    @(621, 'Tuesday, May 18, 2021', 'Arizona Diamondbacks', 0, 'Los Angeles Dodgers', 1),
    @(622, 'Tuesday, May 18, 2021', 'Colorado Rockies', 0, 'San Diego Padres', 0),
    @(623, 'Tuesday, May 18, 2021', 'Detroit Tigers', 1, 'Seattle Mariners', 0)
)

foreach ($row in $textRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 3).Style = "Normal"

    $ws.Cells.Item($r, 4).Value = $row[4]

    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 5).Style = "Normal"
}
